# Update recomputed statistics in the POIs_sorted worksheet.
# Columns: D = GO_POI fraction, G = Value (PPI_rank related score).
# Only the numeric values that actually changed between the previous run
# and the corrected run are touched; every other cell is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 0.015099929233531
$ws.Range("D3").Value  = 0.0540540540540541
$ws.Range("G3").Value  = 5.90806584989691
$ws.Range("D4").Value  = 0.0526315789473684
$ws.Range("G4").Value  = 2.81192265581059
$ws.Range("D5").Value  = 0.0476190476190476
$ws.Range("G5").Value  = 3.20813890652545
$ws.Range("G6").Value  = 2.9381759817832
$ws.Range("G7").Value  = 2.26966306455418
$ws.Range("G8").Value  = 2.23384916725504
$ws.Range("G9").Value  = 1.36020209880448
$ws.Range("G10").Value = 0.332138464774045
$ws.Range("G11").Value = 0.308404536242237
